# Auto-generated edit script applying the Ultima_Profits profit-recalculation update
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets.
$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 113198.336
$ws.Range("I99").Value = 181.33333
$ws.Range("J99").Value = 169706.83
$ws.Range("K99").Value = 543.99999
$ws.Range("L99").Value = 509120.49
$ws.Range("M99").Value = 954.00001
$ws.Range("N99").Value = -512116.49
$ws.Range("H127").Value = 792.73914
$ws.Range("I127").Value = 548.8333
$ws.Range("J127").Value = 878.82355
$ws.Range("K127").Value = 1646.4999
$ws.Range("L127").Value = 2636.47065
$ws.Range("M127").Value = 3313.5001
$ws.Range("N127").Value = -12556.47065
$ws.Range("H131").Value = 11147.143
$ws.Range("I131").Value = 867.8570999999999
$ws.Range("J131").Value = 21426.428
$ws.Range("K131").Value = 2603.5713
$ws.Range("L131").Value = 64279.284
$ws.Range("M131").Value = 2436.4287
$ws.Range("N131").Value = -74359.284
$ws.Range("H137").Value = 8000900.5
$ws.Range("I137").Value = 850.53845
$ws.Range("K137").Value = 2551.61535
$ws.Range("M137").Value = -1.615350000000035
$ws.Range("H138").Value = 1680.4642
$ws.Range("I138").Value = 997.95654
$ws.Range("J138").Value = 4820
$ws.Range("K138").Value = 2993.86962
$ws.Range("L138").Value = 14460
$ws.Range("M138").Value = 2146.13038
$ws.Range("N138").Value = -24740

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 4000
$ws.Range("I64").Value = 4000
$ws.Range("K64").Value = 4000
$ws.Range("M64").Value = -3752
$ws.Range("H67").Value = 4000
$ws.Range("I67").Value = 4000
$ws.Range("K67").Value = 4000
$ws.Range("M67").Value = -3142
$ws.Range("H102").Value = 1464
$ws.Range("I102").Value = 1464
$ws.Range("K102").Value = 1464
$ws.Range("M102").Value = 158

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 45181
$ws.Range("J62").Value = 45181
$ws.Range("L62").Value = 45181
$ws.Range("N62").Value = -46553
$ws.Range("H65").Value = 45181
$ws.Range("J65").Value = 45181
$ws.Range("L65").Value = 135543
$ws.Range("N65").Value = -142407
$ws.Range("H99").Value = 1155.8
$ws.Range("I99").Value = 1162
$ws.Range("J99").Value = 1100
$ws.Range("K99").Value = 1162
$ws.Range("L99").Value = 1100
$ws.Range("M99").Value = 336
$ws.Range("N99").Value = -4096

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7096357.5
$ws.Range("I31").Value = 4994.0625
$ws.Range("J31").Value = 22224600
$ws.Range("K31").Value = 4994.0625
$ws.Range("L31").Value = 22224600
$ws.Range("M31").Value = -4699.0625
$ws.Range("N31").Value = -22225190
$ws.Range("H34").Value = 7096357.5
$ws.Range("I34").Value = 4994.0625
$ws.Range("J34").Value = 22224600
$ws.Range("K34").Value = 4994.0625
$ws.Range("L34").Value = 22224600
$ws.Range("M34").Value = -4792.0625
$ws.Range("N34").Value = -22225004
$ws.Range("H132").Value = 13160179
$ws.Range("I132").Value = 21741216
$ws.Range("J132").Value = 2587.8667
$ws.Range("K132").Value = 65223648
$ws.Range("L132").Value = 7763.6001
$ws.Range("M132").Value = -65221118
$ws.Range("N132").Value = -12823.6001
$ws.Range("H140").Value = 47143.168
$ws.Range("J140").Value = 47143.168
$ws.Range("L140").Value = 47143.168
$ws.Range("N140").Value = -57503.168

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5447.1787
$ws.Range("I3").Value = 3343.6667
$ws.Range("K3").Value = 10031.0001
$ws.Range("M3").Value = -9919.000100000001
$ws.Range("H98").Value = 842.45
$ws.Range("I98").Value = 441.42856
$ws.Range("J98").Value = 1058.3846
$ws.Range("K98").Value = 1324.28568
$ws.Range("L98").Value = 3175.1538
$ws.Range("M98").Value = 173.71432
$ws.Range("N98").Value = -6171.1538
$ws.Range("H129").Value = 2989.68
$ws.Range("I129").Value = 2916.6667
$ws.Range("J129").Value = 2999.6365
$ws.Range("K129").Value = 8750.000100000001
$ws.Range("L129").Value = 8998.9095
$ws.Range("M129").Value = -3750.000100000001
$ws.Range("N129").Value = -18998.9095
$ws.Range("H131").Value = 858.38
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 858.38
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2575.14
$ws.Range("M131").ClearContents()  # was 3309.9999
$ws.Range("N131").Value = -12655.14

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 13333.333
$ws.Range("I64").Value = 10000
$ws.Range("J64").Value = 15000
$ws.Range("K64").Value = 10000
$ws.Range("L64").Value = 15000
$ws.Range("M64").Value = -9752
$ws.Range("N64").Value = -15496
$ws.Range("H67").Value = 13333.333
$ws.Range("I67").Value = 10000
$ws.Range("J67").Value = 15000
$ws.Range("K67").Value = 10000
$ws.Range("L67").Value = 15000
$ws.Range("M67").Value = -9142
$ws.Range("N67").Value = -16716
$ws.Range("H113").Value = 84906
$ws.Range("I113").Value = 84906
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 84906
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -82736
$ws.Range("N113").ClearContents()  # was -5597.5
$ws.Range("H125").Value = 56884
$ws.Range("J125").Value = 56884
$ws.Range("L125").Value = 56884
$ws.Range("N125").Value = -61804
$ws.Range("H128").Value = 52316.668
$ws.Range("J128").Value = 52316.668
$ws.Range("L128").Value = 52316.668
$ws.Range("N128").Value = -62276.668

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7731.5386
$ws.Range("J40").Value = 4319.091
$ws.Range("L40").Value = 4319.091
$ws.Range("N40").Value = -4591.091
$ws.Range("H46").Value = 819.3570999999999
$ws.Range("I46").Value = 760.2
$ws.Range("K46").Value = 760.2
$ws.Range("M46").Value = -572.2
$ws.Range("H122").Value = 4575.8213
$ws.Range("I122").Value = 4895.722
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 14687.166
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -12237.166
$ws.Range("N122").Value = -16900

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 850.4706
$ws.Range("I107").Value = 892.375
$ws.Range("J107").Value = 180
$ws.Range("K107").Value = 2677.125
$ws.Range("L107").Value = 540
$ws.Range("M107").Value = -757.125
$ws.Range("N107").Value = -4380
